# "Generate Report for Handback"
# Populates the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns on the per-locale sheets, and flips the overall status from
# "In Translation" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "In Translation" -> "Handed back: in sync with en-US"
#    (shows up on Overview!E:F and on the locale sheets' Status column C)
# ---------------------------------------------------------------------------
$wsOverview.Cells.Replace("In Translation", "Handed back: in sync with en-US")
$wsZhCn.Cells.Replace("In Translation", "Handed back: in sync with en-US")
$wsDeDe.Cells.Replace("In Translation", "Handed back: in sync with en-US")

# ---------------------------------------------------------------------------
# 2. Widen columns to fit the new handback file names / dates
# ---------------------------------------------------------------------------
# Overview sheet: columns E (zh-cn) and F (de-de)
$wsOverview.Columns.Item(5).ColumnWidth = 29.1111111111111
$wsOverview.Columns.Item(6).ColumnWidth = 29.1111111111111

# Locale sheets: column C (Status), I (Latest Target File), J (Latest Handback File)
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Columns.Item(3).ColumnWidth  = 29.1111111111111
    $ws.Columns.Item(9).ColumnWidth  = 39.1666666666667
    $ws.Columns.Item(10).ColumnWidth = 39.1666666666667
}

# ---------------------------------------------------------------------------
# 3. Latest Handback DateTime (column K)
#    zh-cn handback happened first, de-de a few seconds later, so they end
#    up with two different timestamps even though both started as the same
#    placeholder value.
# ---------------------------------------------------------------------------
$wsZhCn.Cells.Replace("0001-01-01 00:00:00", "2016-08-22 10:04:31")
$wsDeDe.Cells.Replace("0001-01-01 00:00:00", "2016-08-22 10:04:31")
$wsDeDe.Range("K2").Value = "2016-08-22 10:04:38"
$wsDeDe.Range("K3").Value = "2016-08-22 10:04:38"

# ---------------------------------------------------------------------------
# 4. Fill in "Latest Target File" (I) and "Latest Handback File" (J), and
#    turn "Latest Target File" into a hyperlink back to the source .md file
#    (same as column A), for both rows on both locale sheets.
# ---------------------------------------------------------------------------
$mdName1  = "15739a83-f072-4a82-b462-dd07dea86eb3.md"
$mdName2  = "694ee841-4510-43b3-8344-907b04704a1c.md"
$mdUrl1   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a46125dd9cecfedd3c8a9196c1300162115b949/e2e/15739a83-f072-4a82-b462-dd07dea86eb3.md"
$mdUrl2   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a46125dd9cecfedd3c8a9196c1300162115b949/e2e/694ee841-4510-43b3-8344-907b04704a1c.md"

$zhXlf1 = "15739a83-f072-4a82-b462-dd07dea86eb3.8c15a392a039b9333e4a6cfbef261a9677492310.zh-cn.xlf"
$zhXlf2 = "694ee841-4510-43b3-8344-907b04704a1c.f097a3cb591c1afa65bc67fea039831cbfcae728.zh-cn.xlf"
$deXlf1 = "15739a83-f072-4a82-b462-dd07dea86eb3.8c15a392a039b9333e4a6cfbef261a9677492310.de-de.xlf"
$deXlf2 = "694ee841-4510-43b3-8344-907b04704a1c.f097a3cb591c1afa65bc67fea039831cbfcae728.de-de.xlf"

function Set-HandbackRow {
    param($ws, $row, $mdName, $mdUrl, $xlfName)

    $iCell = $ws.Cells.Item($row, 9)
    $jCell = $ws.Cells.Item($row, 10)

    $jCell.Value = $xlfName

    $iCell.Value = $mdName
    $iCell.Font.Underline = 2
    $iCell.Font.Color = 15570276
}

# Re-create all the hyperlinks on each locale sheet so that the "Latest
# Target File" column (I) gets its own hyperlink, interleaved with the
# existing "Source File Name" (A) links in row order (A2, I2, A3, I3).
foreach ($item in @(
        @{ ws = $wsZhCn; xlf1 = $zhXlf1; xlf2 = $zhXlf2 },
        @{ ws = $wsDeDe; xlf1 = $deXlf1; xlf2 = $deXlf2 }
    )) {

    $ws = $item.ws

    Set-HandbackRow -ws $ws -row 2 -mdName $mdName1 -mdUrl $mdUrl1 -xlfName $item.xlf1
    Set-HandbackRow -ws $ws -row 3 -mdName $mdName2 -mdUrl $mdUrl2 -xlfName $item.xlf2

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $mdName1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl1, "", "", $mdName1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $mdName2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdUrl2, "", "", $mdName2)
}
